# Update TPM-derived specificity values in the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 10.45366133333333
$ws.Range("N2").Value = 31.360984
$ws.Range("O2").Value = 0.2775334303572506
$ws.Range("P2").Value = 0.2775334303572506
$ws.Range("Q2").Value = 1.211426028167111
$ws.Range("R2").Value = 10.902834253504
$ws.Range("S2").Value = 0.2775334303572506
$ws.Range("T2").Value = 0.2775334303572506

# Row 3 (only specificity columns changed)
$ws.Range("O3").Value = 0.2992222971432776
$ws.Range("P3").Value = 0.2992222971432776
$ws.Range("S3").Value = 0.2992222971432776
$ws.Range("T3").Value = 0.2992222971432776

# Row 4
$ws.Range("M4").Value = 6.998235333333334
$ws.Range("N4").Value = 20.994706
$ws.Range("O4").Value = 0.1857955979800236
$ws.Range("P4").Value = 0.1857955979800235
$ws.Range("Q4").Value = 0.8109928343484445
$ws.Range("R4").Value = 7.298935509136001
$ws.Range("S4").Value = 0.1857955979800236
$ws.Range("T4").Value = 0.1857955979800235

# Row 5
$ws.Range("M5").Value = 3.843654
$ws.Range("N5").Value = 11.530962
$ws.Range("O5").Value = 0.1020448669333559
$ws.Range("P5").Value = 0.1020448669333559
$ws.Range("Q5").Value = 0.4454231250080001
$ws.Range("R5").Value = 4.008808125072
$ws.Range("S5").Value = 0.1020448669333559
$ws.Range("T5").Value = 0.1020448669333559

# Row 6
$ws.Range("M6").Value = 5.100162333333333
$ws.Range("N6").Value = 15.300487
$ws.Range("O6").Value = 0.1354038075860923
$ws.Range("P6").Value = 0.1354038075860922
$ws.Range("Q6").Value = 0.5910340120524445
$ws.Range("R6").Value = 5.319306108472
$ws.Range("S6").Value = 0.1354038075860923
$ws.Range("T6").Value = 0.1354038075860922
